{"js": "const body = context.document.body;\n\nconst replacements = [\n    [\"2025-04-09 Wednesday\", \"2025-04-10 Thursday\"],\n    [\"887\u00d79=7983\", \"927\u00d78=7416\"],\n    [\"757\u00d77=5299\", \"558\u00d78=4464\"],\n    [\"912\u00d72=1824\", \"278\u00d79=2502\"],\n    [\"705\u00d76=4230\", \"168\u00d79=1512\"],\n    [\"705\u00d73=2115\", \"328\u00d78=2624\"],\n    [\"701\u00d75=3505\", \"249\u00d78=1992\"],\n    [\"789\u00d73=2367\", \"659\u00d79=5931\"],\n    [\"792\u00d73=2376\", \"401\u00d74=1604\"],\n    [\"897\u00d76=5382\", \"119\u00d74=476\"],\n    [\"880\u00d75=4400\", \"231\u00d72=462\"],\n    [\"201\u00d75=1005\", \"574\u00d74=2296\"],\n    [\"850\u00d78=6800\", \"628\u00d76=3768\"],\n    [\"695\u00d72=1390\", \"444\u00d76=2664\"],\n    [\"462\u00d75=2310\", \"935\u00d74=3740\"],\n    [\"290\u00d78=2320\", \"417\u00d76=2502\"],\n    [\"936\u00d74=3744\", \"954\u00d77=6678\"],\n    [\"531\u00d74=2124\", \"790\u00d72=1580\"],\n    [\"887\u00d76=5322\", \"175\u00d77=1225\"],\n    [\"713\u00d72=1426\", \"480\u00d73=1440\"],\n    [\"387\u00d75=1935\", \"849\u00d75=4245\"],\n    [\"959\u00d79=8631\", \"106\u00d74=424\"],\n    [\"360\u00d72=720\", \"671\u00d79=6039\"],\n    [\"390\u00d76=2340\", \"671\u00d74=2684\"],\n    [\"456\u00d75=2280\", \"840\u00d79=7560\"],\n    [\"963\u00d76=5778\", \"683\u00d76=4098\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n    for (const item of results.items) {\n        item.insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2025-04-09 Wednesday\"; New = \"2025-04-10 Thursday\"},\n    @{Old = \"887\u00d79=7983\"; New = \"927\u00d78=7416\"},\n    @{Old = \"757\u00d77=5299\"; New = \"558\u00d78=4464\"},\n    @{Old = \"912\u00d72=1824\"; New = \"278\u00d79=2502\"},\n    @{Old = \"705\u00d76=4230\"; New = \"168\u00d79=1512\"},\n    @{Old = \"705\u00d73=2115\"; New = \"328\u00d78=2624\"},\n    @{Old = \"701\u00d75=3505\"; New = \"249\u00d78=1992\"},\n    @{Old = \"789\u00d73=2367\"; New = \"659\u00d79=5931\"},\n    @{Old = \"792\u00d73=2376\"; New = \"401\u00d74=1604\"},\n    @{Old = \"897\u00d76=5382\"; New = \"119\u00d74=476\"},\n    @{Old = \"880\u00d75=4400\"; New = \"231\u00d72=462\"},\n    @{Old = \"201\u00d75=1005\"; New = \"574\u00d74=2296\"},\n    @{Old = \"850\u00d78=6800\"; New = \"628\u00d76=3768\"},\n    @{Old = \"695\u00d72=1390\"; New = \"444\u00d76=2664\"},\n    @{Old = \"462\u00d75=2310\"; New = \"935\u00d74=3740\"},\n    @{Old = \"290\u00d78=2320\"; New = \"417\u00d76=2502\"},\n    @{Old = \"936\u00d74=3744\"; New = \"954\u00d77=6678\"},\n    @{Old = \"531\u00d74=2124\"; New = \"790\u00d72=1580\"},\n    @{Old = \"887\u00d76=5322\"; New = \"175\u00d77=1225\"},\n    @{Old = \"713\u00d72=1426\"; New = \"480\u00d73=1440\"},\n    @{Old = \"387\u00d75=1935\"; New = \"849\u00d75=4245\"},\n    @{Old = \"959\u00d79=8631\"; New = \"106\u00d74=424\"},\n    @{Old = \"360\u00d72=720\"; New = \"671\u00d79=6039\"},\n    @{Old = \"390\u00d76=2340\"; New = \"671\u00d74=2684\"},\n    @{Old = \"456\u00d75=2280\"; New = \"840\u00d79=7560\"},\n    @{Old = \"963\u00d76=5778\"; New = \"683\u00d76=4098\"}\n)\n\nforeach ($rep in $replacements) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Execute($rep.Old, $false, $true, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n}\n"}
